$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Login" (sheet1)
# ---------------------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("Login")

# G12 gains the "HLSupervisorContraCosta" value (previously blank).
$wsLogin.Range("G12").Value = "HLSupervisorContraCosta"

# B13 iteration count bumped from 1 to 2.
$wsLogin.Range("B13").Value = 2

# New test-case row appended at the bottom of the table (row 21).
$wsLogin.Range("A21").Value = "testT4250"
$wsLogin.Range("B21").Value = 1
$wsLogin.Range("C21").Value = 1
$wsLogin.Range("F21").Value = "HLWorkerContraCosta"
$wsLogin.Range("G21").Value = "HLSupervisorContraCosta"
$wsLogin.Range("H21").Value = "ERWorkerContraCosta"

# Match the formatting used by the other rows of the same table so the new
# row 21 looks consistent with its neighbours.
$wsLogin.Range("G5").Copy()
$wsLogin.Range("A21").PasteSpecial(-4122)
$wsLogin.Range("A5").Copy()
$wsLogin.Range("B21:C21").PasteSpecial(-4122)
$wsLogin.Range("G21").PasteSpecial(-4122)
$wsLogin.Range("E20").Copy()
$wsLogin.Range("D21:E21").PasteSpecial(-4122)
$wsLogin.Range("F21").PasteSpecial(-4122)
$wsLogin.Range("H21:I21").PasteSpecial(-4122)
$wsLogin.Application.CutCopyMode = $false

# Rows in this table use a 15pt custom row height; match it for the new row.
$wsLogin.Rows.Item(21).RowHeight = 15

# ---------------------------------------------------------------------------
# Sheet "Logout" (sheet2)
# ---------------------------------------------------------------------------
$wsLogout = $wb.Worksheets.Item("Logout")

# New test-case row appended at the bottom of the table (row 8).
$wsLogout.Range("A8").Value = "testT4250"
$wsLogout.Range("B8").Value = 1
$wsLogout.Range("C8").Value = 1
$wsLogout.Range("D8").Value = "Click"
$wsLogout.Range("E8").Value = "Click"
$wsLogout.Range("F8").Value = "Yes"

# Formatting consistent with the rest of the table (rows 1-7 all use the
# same style).
$wsLogout.Range("A7").Copy()
$wsLogout.Range("A8:F8").PasteSpecial(-4122)
$wsLogout.Application.CutCopyMode = $false

[void]$wsLogout.Range("F8").Select()

# Re-activate the "Login" tab (it was the selected tab before the edits) and
# restore its own last-used selection.
[void]$wsLogin.Range("D9").Select()
